# The "mit ettem" (what did I eat) tracker sheet had a header row (row 1)
# meant to describe repeating "portion, food" column pairs (matching the
# data columns actually used below it: C/D, E/F, G/H, I/J, K/L).
# In the original file the header was inconsistent: columns E and G were
# left blank, columns I and K incorrectly said "food" instead of "portion",
# and the pattern needlessly continued (with "food") out to column P.
# Fix the header row so that:
#   - the previously-empty E1 and G1 cells get the "portion" label
#   - I1 and K1 switch from "food" to "portion" (completing the pattern)
#   - the now-unused trailing M1:P1 header cells are removed entirely
# This shrinks the sheet's used range from A1:P13 down to A1:L13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing header labels / fix the alternating portion-food pattern.
$ws.Range("E1").Value = "portion"
$ws.Range("G1").Value = "portion"
$ws.Range("I1").Value = "portion"
$ws.Range("K1").Value = "portion"

# Remove the now-superfluous trailing header cells (M1:P1) completely,
# shrinking the sheet's dimension/used range to A1:L13.
$ws.Range("M1:P1").Clear()

# The author's last selection before saving ended up on G23.
$ws.Range("G23").Select()
